$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("M3").Value = "Commande Attente"
$ws.Range("O3").ClearContents()

# Row 5
$ws.Range("M5").Value = "Client"
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("G6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("M6").Value = "Commande"
$ws.Range("O6").ClearContents()

# Row 8
$ws.Range("E8").Value = "Client"

# Row 9
$ws.Range("E9").ClearContents()
$ws.Range("I9").Value = "Garniture_Commande"

# Row 10
$ws.Range("I10").Value = "Commande"

# Row 11
$ws.Range("I11").Value = "Garniture"

# View adjustments
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("I12").Select()
